# Auto-generated edit script applying the Atomos_Profits.xlsx diff
# Updates market-price-derived columns (H,I,J,K,L,M,N) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3069.25
$ws.Range("J64").Value = 2983.3333
$ws.Range("L64").Value = 2983.3333
$ws.Range("N64").Value = -3479.3333
$ws.Range("H67").Value = 3069.25
$ws.Range("J67").Value = 2983.3333
$ws.Range("L67").Value = 2983.3333
$ws.Range("N67").Value = -4699.3333
$ws.Range("H74").Value = 3954.3125
$ws.Range("I74").Value = 3407.25
$ws.Range("J74").Value = 4136.6665
$ws.Range("K74").Value = 3407.25
$ws.Range("L74").Value = 4136.6665
$ws.Range("M74").Value = -2471.25
$ws.Range("N74").Value = -6008.6665
$ws.Range("H76").Value = 3000.4
$ws.Range("I76").Value = 2799.875
$ws.Range("K76").Value = 2799.875
$ws.Range("M76").Value = -2484.875
$ws.Range("H77").Value = 3954.3125
$ws.Range("I77").Value = 3407.25
$ws.Range("J77").Value = 4136.6665
$ws.Range("K77").Value = 17036.25
$ws.Range("L77").Value = 20683.3325
$ws.Range("M77").Value = -12356.25
$ws.Range("N77").Value = -30043.3325
$ws.Range("H79").Value = 3000.4
$ws.Range("I79").Value = 2799.875
$ws.Range("K79").Value = 2799.875
$ws.Range("M79").Value = -1707.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1476.305
$ws.Range("I45").Value = 1130.5366
$ws.Range("K45").Value = 1130.5366
$ws.Range("M45").Value = -753.5365999999999
$ws.Range("H63").Value = 2494.1177
$ws.Range("I63").Value = 2125
$ws.Range("J63").Value = 3380
$ws.Range("K63").Value = 2125
$ws.Range("L63").Value = 3380
$ws.Range("M63").Value = -1439
$ws.Range("N63").Value = -4752
$ws.Range("H66").Value = 2494.1177
$ws.Range("I66").Value = 2125
$ws.Range("J66").Value = 3380
$ws.Range("K66").Value = 10625
$ws.Range("L66").Value = 16900
$ws.Range("M66").Value = -7193
$ws.Range("N66").Value = -23764
$ws.Range("H80").Value = 24860.777
$ws.Range("J80").Value = 26705.875
$ws.Range("L80").Value = 26705.875
$ws.Range("N80").Value = -28701.875
$ws.Range("H83").Value = 24860.777
$ws.Range("J83").Value = 26705.875
$ws.Range("L83").Value = 80117.625
$ws.Range("N83").Value = -90101.625
$ws.Range("H86").Value = 40000
$ws.Range("J86").Value = 40000
$ws.Range("L86").Value = 40000
$ws.Range("N86").Value = -42372
$ws.Range("H88").Value = 1519.2
$ws.Range("I88").Value = 1519.2
$ws.Range("K88").Value = 1519.2
$ws.Range("M88").Value = -1113.2
$ws.Range("H89").Value = 40000
$ws.Range("J89").Value = 40000
$ws.Range("L89").Value = 120000
$ws.Range("N89").Value = -131856
$ws.Range("H91").Value = 1519.2
$ws.Range("I91").Value = 1519.2
$ws.Range("K91").Value = 1519.2
$ws.Range("M91").Value = -115.2
$ws.Range("H140").Value = 29642.857
$ws.Range("J140").Value = 29642.857
$ws.Range("L140").Value = 29642.857
$ws.Range("N140").Value = -40002.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H82").Value = 19868.588
$ws.Range("I82").Value = 4592.3335
$ws.Range("J82").Value = 28201.092
$ws.Range("K82").Value = 4592.3335
$ws.Range("L82").Value = 28201.092
$ws.Range("M82").Value = -4209.3335
$ws.Range("N82").Value = -28967.092
$ws.Range("H85").Value = 19868.588
$ws.Range("I85").Value = 4592.3335
$ws.Range("J85").Value = 28201.092
$ws.Range("K85").Value = 4592.3335
$ws.Range("L85").Value = 28201.092
$ws.Range("M85").Value = -3266.3335
$ws.Range("N85").Value = -30853.092
$ws.Range("H86").Value = 28053.5
$ws.Range("I86").Value = 2610.5454
$ws.Range("K86").Value = 2610.5454
$ws.Range("M86").Value = -1487.5454
$ws.Range("H89").Value = 28053.5
$ws.Range("I89").Value = 2610.5454
$ws.Range("K89").Value = 13052.727
$ws.Range("M89").Value = -7436.726999999999
$ws.Range("H105").Value = 1491.3513
$ws.Range("I105").Value = 1462.7273
$ws.Range("J105").Value = 1533.3334
$ws.Range("K105").Value = 1462.7273
$ws.Range("L105").Value = 1533.3334
$ws.Range("M105").Value = 284.2727
$ws.Range("N105").Value = -5027.3334
$ws.Range("H107").Value = 2522.75
$ws.Range("I107").Value = 900
$ws.Range("K107").Value = 900
$ws.Range("M107").Value = 1020

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4134.8096
$ws.Range("I62").Value = 3176.5625
$ws.Range("K62").Value = 3176.5625
$ws.Range("M62").Value = -2552.5625
$ws.Range("H65").Value = 4134.8096
$ws.Range("I65").Value = 3176.5625
$ws.Range("K65").Value = 15882.8125
$ws.Range("M65").Value = -12762.8125
$ws.Range("H105").Value = 4608.3335
$ws.Range("I105").Value = 4355.5557
$ws.Range("J105").Value = 5366.6665
$ws.Range("K105").Value = 4355.5557
$ws.Range("L105").Value = 5366.6665
$ws.Range("M105").Value = -2608.5557
$ws.Range("N105").Value = -8860.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1039.356
$ws.Range("I131").Value = 551.6667
$ws.Range("J131").Value = 1094.566
$ws.Range("K131").Value = 1655.0001
$ws.Range("L131").Value = 3283.698
$ws.Range("M131").Value = 3384.9999
$ws.Range("N131").Value = -13363.698

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1334.0426
$ws.Range("I46").Value = 970
$ws.Range("K46").Value = 970
$ws.Range("M46").Value = -782
$ws.Range("H68").Value = 1811.3043
$ws.Range("I68").Value = 1043
$ws.Range("K68").Value = 1043
$ws.Range("M68").Value = -294
$ws.Range("H71").Value = 1811.3043
$ws.Range("I71").Value = 1043
$ws.Range("K71").Value = 5215
$ws.Range("M71").Value = -1471
$ws.Range("H82").Value = 3616
$ws.Range("I82").Value = 2452
$ws.Range("K82").Value = 2452
$ws.Range("M82").Value = -2091
$ws.Range("H85").Value = 3616
$ws.Range("I85").Value = 2452
$ws.Range("K85").Value = 2452
$ws.Range("M85").Value = -1204

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 27095.666
$ws.Range("J124").Value = 27095.666
$ws.Range("L124").Value = 27095.666
$ws.Range("N124").Value = -36915.666
